$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "17586324521"
$ws.Range("D2").Value = "Cloris629k"
$ws.Range("E2").Value = "Alicerlzq"
